# Automatische test-sync: 2025-08-03 15:00:50
#
# Appends the 13th test-mail row to the "Logs" sheet, appends the matching
# category tally row to the "Dashboard" sheet, widens the conditional
# formatting ranges to cover the new row, and extends the bar chart's
# category/value series so it picks up the new "Documentatie / Datasheets"
# bucket.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 21
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A21").Value = "Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Range("B21").Value = "mailmind.test@zohomail.eu"
$logs.Range("C21").Value = "Testmail #13: Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Range("D21").Value = "Documentatie / Datasheets"
$logs.Range("E21").Value = "Bedankt, we hebben dit doorgestuurd naar documentatie@bedrijf.nl."
$logs.Range("F21").Value = "2025-08-03 14:59:51"
$logs.Range("G21").Value = "Ja"
$logs.Range("H21").Value = "Ja"
$logs.Range("I21").Value = "Nee"
$logs.Range("J21").Value = "Nee"

# Widen the conditional-formatting blocks (D/G/H/I/J, rows 2-20 -> 2-21) so
# the new row inherits the same category / yes-no colouring rules.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "20")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "21")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append row 6 (new category tally)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Documentatie / Datasheets"
$dash.Range("B6").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend category/value series from row 5 to row 6
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "=Dashboard!`$A`$2:`$A`$6"
$series.Values = "=Dashboard!`$B`$2:`$B`$6"
